$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(412, 1).Value = 'A-B'
$ws.Cells.Item(412, 2).Value = 1
$ws.Cells.Item(412, 3).Value = 'Unambiguous'
$ws.Cells.Item(412, 4).Value = 2.3488888888888888
$ws.Cells.Item(412, 5).Value = 10.31111111111111
$ws.Cells.Item(412, 6).Value = 'A'
$ws.Cells.Item(412, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(412, 9).Value = 'Current Work'
$ws.Cells.Item(412, 10).Value = 1

$ws.Cells.Item(413, 1).Value = 'A-B'
$ws.Cells.Item(413, 2).Value = 2
$ws.Cells.Item(413, 3).Value = 'Unambiguous'
$ws.Cells.Item(413, 4).Value = -1.753125
$ws.Cells.Item(413, 5).Value = 13.153124999999999
$ws.Cells.Item(413, 6).Value = 'B'
$ws.Cells.Item(413, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(413, 9).Value = 'Current Work'
$ws.Cells.Item(413, 10).Value = 1

$ws.Cells.Item(414, 1).Value = 'B-H'
$ws.Cells.Item(414, 2).Value = 1
$ws.Cells.Item(414, 3).Value = 'Unambiguous'
$ws.Cells.Item(414, 4).Value = -1.753125
$ws.Cells.Item(414, 5).Value = 13.153124999999999
$ws.Cells.Item(414, 6).Value = 'B'
$ws.Cells.Item(414, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(414, 9).Value = 'Current Work'
$ws.Cells.Item(414, 10).Value = 1

$ws.Cells.Item(415, 1).Value = 'B-H'
$ws.Cells.Item(415, 2).Value = 2
$ws.Cells.Item(415, 3).Value = 'Unambiguous'
$ws.Cells.Item(415, 4).Value = -5.5714285714285703
$ws.Cells.Item(415, 5).Value = 15.75
$ws.Cells.Item(415, 6).Value = 'H'
$ws.Cells.Item(415, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(415, 9).Value = 'Current Work'
$ws.Cells.Item(415, 10).Value = 1

$ws.Cells.Item(416, 1).Value = 'B-C'
$ws.Cells.Item(416, 2).Value = 1
$ws.Cells.Item(416, 3).Value = 'Unambiguous'
$ws.Cells.Item(416, 4).Value = -1.753125
$ws.Cells.Item(416, 5).Value = 13.153124999999999
$ws.Cells.Item(416, 6).Value = 'B'
$ws.Cells.Item(416, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(416, 9).Value = 'Current Work'
$ws.Cells.Item(416, 10).Value = 1

$ws.Cells.Item(417, 1).Value = 'B-C'
$ws.Cells.Item(417, 2).Value = 2
$ws.Cells.Item(417, 3).Value = 'Unambiguous'
$ws.Cells.Item(417, 4).Value = -0.40769230769230752
$ws.Cells.Item(417, 5).Value = 21.388461538461538
$ws.Cells.Item(417, 6).Value = 'C'
$ws.Cells.Item(417, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(417, 9).Value = 'Current Work'
$ws.Cells.Item(417, 10).Value = 1

$ws.Cells.Item(418, 1).Value = 'C-D'
$ws.Cells.Item(418, 2).Value = 1
$ws.Cells.Item(418, 3).Value = 'Ambiguous'
$ws.Cells.Item(418, 4).Value = -0.40769230769230752
$ws.Cells.Item(418, 5).Value = 21.388461538461538
$ws.Cells.Item(418, 6).Value = 'C'
$ws.Cells.Item(418, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(418, 9).Value = 'Current Work'
$ws.Cells.Item(418, 10).Value = 1

$ws.Cells.Item(419, 1).Value = 'C-D'
$ws.Cells.Item(419, 2).Value = 2
$ws.Cells.Item(419, 3).Value = 'Ambiguous'
$ws.Cells.Item(419, 4).Value = -1.78
$ws.Cells.Item(419, 5).Value = 27.603999999999996
$ws.Cells.Item(419, 6).Value = 'D'
$ws.Cells.Item(419, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(419, 9).Value = 'Current Work'
$ws.Cells.Item(419, 10).Value = 1

$ws.Cells.Item(420, 1).Value = 'D-E'
$ws.Cells.Item(420, 2).Value = 1
$ws.Cells.Item(420, 3).Value = 'Ambiguous'
$ws.Cells.Item(420, 4).Value = -1.78
$ws.Cells.Item(420, 5).Value = 27.603999999999996
$ws.Cells.Item(420, 6).Value = 'D'
$ws.Cells.Item(420, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(420, 9).Value = 'Current Work'
$ws.Cells.Item(420, 10).Value = 1

$ws.Cells.Item(421, 1).Value = 'D-E'
$ws.Cells.Item(421, 2).Value = 2
$ws.Cells.Item(421, 3).Value = 'Ambiguous'
$ws.Cells.Item(421, 4).Value = -2.35
$ws.Cells.Item(421, 5).Value = 38.449999999999996
$ws.Cells.Item(421, 6).Value = 'E'
$ws.Cells.Item(421, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(421, 9).Value = 'Current Work'
$ws.Cells.Item(421, 10).Value = 1

$ws.Cells.Item(422, 1).Value = 'E-G'
$ws.Cells.Item(422, 2).Value = 1
$ws.Cells.Item(422, 3).Value = 'Ambiguous'
$ws.Cells.Item(422, 4).Value = -1.78
$ws.Cells.Item(422, 5).Value = 27.603999999999996
$ws.Cells.Item(422, 6).Value = 'D'
$ws.Cells.Item(422, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(422, 9).Value = 'Current Work'
$ws.Cells.Item(422, 10).Value = 1

$ws.Cells.Item(423, 1).Value = 'E-G'
$ws.Cells.Item(423, 2).Value = 2
$ws.Cells.Item(423, 3).Value = 'Ambiguous'
$ws.Cells.Item(423, 4).Value = -5.85
$ws.Cells.Item(423, 5).Value = 37.712499999999999
$ws.Cells.Item(423, 6).Value = 'G'
$ws.Cells.Item(423, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(423, 9).Value = 'Current Work'
$ws.Cells.Item(423, 10).Value = 1

$ws.Cells.Item(424, 1).Value = 'E-F'
$ws.Cells.Item(424, 2).Value = 1
$ws.Cells.Item(424, 3).Value = 'Ambiguous'
$ws.Cells.Item(424, 4).Value = -2.35
$ws.Cells.Item(424, 5).Value = 38.449999999999996
$ws.Cells.Item(424, 6).Value = 'E'
$ws.Cells.Item(424, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(424, 9).Value = 'Current Work'
$ws.Cells.Item(424, 10).Value = 1

$ws.Cells.Item(425, 1).Value = 'E-F'
$ws.Cells.Item(425, 2).Value = 2
$ws.Cells.Item(425, 3).Value = 'Ambiguous'
$ws.Cells.Item(425, 4).Value = -3.7333333333333329
$ws.Cells.Item(425, 5).Value = 33
$ws.Cells.Item(425, 6).Value = 'F'
$ws.Cells.Item(425, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(425, 9).Value = 'Current Work'
$ws.Cells.Item(425, 10).Value = 1

$ws.Cells.Item(426, 1).Value = 'E-J'
$ws.Cells.Item(426, 2).Value = 1
$ws.Cells.Item(426, 3).Value = 'Ambiguous'
$ws.Cells.Item(426, 4).Value = -2.35
$ws.Cells.Item(426, 5).Value = 38.449999999999996
$ws.Cells.Item(426, 6).Value = 'E'
$ws.Cells.Item(426, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(426, 9).Value = 'Current Work'
$ws.Cells.Item(426, 10).Value = 1

$ws.Cells.Item(427, 1).Value = 'E-J'
$ws.Cells.Item(427, 2).Value = 2
$ws.Cells.Item(427, 3).Value = 'Ambiguous'
$ws.Cells.Item(427, 4).Value = -1.44444444444444
$ws.Cells.Item(427, 5).Value = 30.911111111111115
$ws.Cells.Item(427, 6).Value = 'J'
$ws.Cells.Item(427, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(427, 9).Value = 'Current Work'
$ws.Cells.Item(427, 10).Value = 1

$ws.Cells.Item(428, 1).Value = 'E-L'
$ws.Cells.Item(428, 2).Value = 1
$ws.Cells.Item(428, 3).Value = 'Ambiguous'
$ws.Cells.Item(428, 4).Value = -2.35
$ws.Cells.Item(428, 5).Value = 38.449999999999996
$ws.Cells.Item(428, 6).Value = 'E'
$ws.Cells.Item(428, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(428, 9).Value = 'Current Work'
$ws.Cells.Item(428, 10).Value = 1

$ws.Cells.Item(429, 1).Value = 'E-L'
$ws.Cells.Item(429, 2).Value = 2
$ws.Cells.Item(429, 3).Value = 'Ambiguous'
$ws.Cells.Item(429, 4).Value = -9.16
$ws.Cells.Item(429, 5).Value = 26.439999999999998
$ws.Cells.Item(429, 6).Value = 'L'
$ws.Cells.Item(429, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(429, 9).Value = 'Current Work'
$ws.Cells.Item(429, 10).Value = 1

$ws.Cells.Item(430, 1).Value = 'L-K'
$ws.Cells.Item(430, 2).Value = 1
$ws.Cells.Item(430, 3).Value = 'Unambiguous'
$ws.Cells.Item(430, 4).Value = -9.16
$ws.Cells.Item(430, 5).Value = 26.439999999999998
$ws.Cells.Item(430, 6).Value = 'L'
$ws.Cells.Item(430, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(430, 9).Value = 'Current Work'
$ws.Cells.Item(430, 10).Value = 1

$ws.Cells.Item(431, 1).Value = 'L-K'
$ws.Cells.Item(431, 2).Value = 2
$ws.Cells.Item(431, 3).Value = 'Unambiguous'
$ws.Cells.Item(431, 4).Value = -12.962999999999999
$ws.Cells.Item(431, 5).Value = 22.224999999999998
$ws.Cells.Item(431, 6).Value = 'K'
$ws.Cells.Item(431, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(431, 9).Value = 'Current Work'
$ws.Cells.Item(431, 10).Value = 1

$ws.Cells.Item(432, 1).Value = 'K-R'
$ws.Cells.Item(432, 2).Value = 1
$ws.Cells.Item(432, 3).Value = 'Unambiguous'
$ws.Cells.Item(432, 4).Value = -12.962999999999999
$ws.Cells.Item(432, 5).Value = 22.224999999999998
$ws.Cells.Item(432, 6).Value = 'K'
$ws.Cells.Item(432, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(432, 9).Value = 'Current Work'
$ws.Cells.Item(432, 10).Value = 1

$ws.Cells.Item(433, 1).Value = 'K-R'
$ws.Cells.Item(433, 2).Value = 2
$ws.Cells.Item(433, 3).Value = 'Unambiguous'
$ws.Cells.Item(433, 4).Value = -15.824999999999999
$ws.Cells.Item(433, 5).Value = 15.45
$ws.Cells.Item(433, 6).Value = 'R'
$ws.Cells.Item(433, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(433, 9).Value = 'Current Work'
$ws.Cells.Item(433, 10).Value = 1

$ws.Cells.Item(434, 1).Value = 'L-M'
$ws.Cells.Item(434, 2).Value = 1
$ws.Cells.Item(434, 3).Value = 'Unambiguous'
$ws.Cells.Item(434, 4).Value = -9.16
$ws.Cells.Item(434, 5).Value = 26.439999999999998
$ws.Cells.Item(434, 6).Value = 'L'
$ws.Cells.Item(434, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(434, 9).Value = 'Current Work'
$ws.Cells.Item(434, 10).Value = 1

$ws.Cells.Item(435, 1).Value = 'L-M'
$ws.Cells.Item(435, 2).Value = 2
$ws.Cells.Item(435, 3).Value = 'Unambiguous'
$ws.Cells.Item(435, 4).Value = -13.122222222222224
$ws.Cells.Item(435, 5).Value = 29.1944444444444
$ws.Cells.Item(435, 6).Value = 'M'
$ws.Cells.Item(435, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(435, 9).Value = 'Current Work'
$ws.Cells.Item(435, 10).Value = 1

$ws.Cells.Item(436, 1).Value = 'M-N'
$ws.Cells.Item(436, 2).Value = 1
$ws.Cells.Item(436, 3).Value = 'Unambiguous'
$ws.Cells.Item(436, 4).Value = -13.122222222222224
$ws.Cells.Item(436, 5).Value = 29.1944444444444
$ws.Cells.Item(436, 6).Value = 'M'
$ws.Cells.Item(436, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(436, 9).Value = 'Current Work'
$ws.Cells.Item(436, 10).Value = 1

$ws.Cells.Item(437, 1).Value = 'M-N'
$ws.Cells.Item(437, 2).Value = 2
$ws.Cells.Item(437, 3).Value = 'Unambiguous'
$ws.Cells.Item(437, 4).Value = -13.557142857142859
$ws.Cells.Item(437, 5).Value = 33.514285714285712
$ws.Cells.Item(437, 6).Value = 'N'
$ws.Cells.Item(437, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(437, 9).Value = 'Current Work'
$ws.Cells.Item(437, 10).Value = 1

$ws.Cells.Item(438, 1).Value = 'N-P'
$ws.Cells.Item(438, 2).Value = 1
$ws.Cells.Item(438, 3).Value = 'Unambiguous'
$ws.Cells.Item(438, 4).Value = -13.557142857142859
$ws.Cells.Item(438, 5).Value = 33.514285714285712
$ws.Cells.Item(438, 6).Value = 'N'
$ws.Cells.Item(438, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(438, 9).Value = 'Current Work'
$ws.Cells.Item(438, 10).Value = 1

$ws.Cells.Item(439, 1).Value = 'N-P'
$ws.Cells.Item(439, 2).Value = 2
$ws.Cells.Item(439, 3).Value = 'Unambiguous'
$ws.Cells.Item(439, 4).Value = -13.3
$ws.Cells.Item(439, 5).Value = 39.125
$ws.Cells.Item(439, 6).Value = 'P'
$ws.Cells.Item(439, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(439, 9).Value = 'Current Work'
$ws.Cells.Item(439, 10).Value = 1

$ws.Cells.Item(440, 1).Value = 'M-S'
$ws.Cells.Item(440, 2).Value = 1
$ws.Cells.Item(440, 3).Value = 'Unambiguous'
$ws.Cells.Item(440, 4).Value = -13.122222222222224
$ws.Cells.Item(440, 5).Value = 29.1944444444444
$ws.Cells.Item(440, 6).Value = 'M'
$ws.Cells.Item(440, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(440, 9).Value = 'Current Work'
$ws.Cells.Item(440, 10).Value = 1

$ws.Cells.Item(441, 1).Value = 'M-S'
$ws.Cells.Item(441, 2).Value = 2
$ws.Cells.Item(441, 3).Value = 'Unambiguous'
$ws.Cells.Item(441, 4).Value = -23.125
$ws.Cells.Item(441, 5).Value = 29.665624999999999
$ws.Cells.Item(441, 6).Value = 'S'
$ws.Cells.Item(441, 8).Value = 'AMNH - Combined 400 -11-H'
$ws.Cells.Item(441, 9).Value = 'Current Work'
$ws.Cells.Item(441, 10).Value = 1
